# "Added last minute updates"
#
# The first paragraph of the document is a hidden placeholder/ID line:
#   "**ID__AFFARS_5309_topic_9__ID**" + a trailing space run
# It needs to become:
#   "**ID__AFFARS_5309_206_1__ID**"   (no trailing space run)
# and the paragraph also picks up the same paragraph border + left indent
# already used by the other body paragraphs later in the document.

$d = $word.ActiveDocument

# The paragraph's text is:
#   "**ID__AFFARS_5309_topic_9__ID**" (31 chars, run 1) + " " (1 char, run 2) + pilcrow
# Delete the lone trailing-space run (characters 31-32) first so the later
# Find/Replace leaves a single clean run with no xml:space="preserve" leftover.
$spaceRun = $d.Range(31, 32)
$spaceRun.Delete()

# Swap the placeholder id text for the real topic id.
$d.Content.Find.Execute("**ID__AFFARS_5309_topic_9__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5309_206_1__ID**", 2)

# Give paragraph 1 the same paragraph border + left indent used elsewhere
# in the document (5-twip space border on all 4 sides, 225-twip/11.25pt
# left indent instead of the old 120-twip indent).
$p1 = $d.Paragraphs(1)
$pf = $p1.Range.ParagraphFormat
$pf.LeftIndent = 11.25

$borders = $pf.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
